$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: Acknowledgement / cover page - merge the "Jairaj P & Abhishek J
# ... (1BG16CS045 & 1BG16CS003)" runs (which were split by gramStart/gramEnd
# proofErr markers) into a single run with the identical combined text.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Jairaj*1BG16CS045*") {
        $r = $p.Range
        $r.End = $r.End - 1   # exclude the paragraph mark
        # Two-step set: first to an interim value so the real final text is
        # recognised as an actual change (it was already equal to the
        # concatenation of the existing runs), forcing the engine to rebuild
        # the paragraph as a single merged run and drop the proofErr markers.
        $r.Text = "Jairaj P & Abhishek J                                                                                       (1BG16CS045 & 1BG16CS003)`u{2060}"
        $r2 = $p.Range
        $r2.End = $r2.End - 1
        $r2.Text = "Jairaj P & Abhishek J                                                                                       (1BG16CS045 & 1BG16CS003)"
    }
}

# ---------------------------------------------------------------------
# Change 2: remove the stray _GoBack bookmark sitting in front of the
# "CONTENTS" heading.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# Change 3: fix the "REQUIREMNTS" typo -> "REQUIREMENTS" in the table of
# contents entry "2.1 HARDWARE AND SOFTWARE REQUIREMNTS ... 8", and leave
# the _GoBack bookmark positioned right after the newly typed "E" (mirroring
# Word's habit of dropping _GoBack at the last edited location).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*HARDWARE AND SOFTWARE REQUIREMNTS*") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "2.1 HARDWARE AND SOFTWARE REQUIREMENTS                                        8 "
    }
}

$rng = $d.Content
$rng.Find.Execute("2.1 HARDWARE AND SOFTWARE REQUIREMENTS") | Out-Null
$start = $rng.Start

# Split off "E" from the rest of the run using a throwaway bookmark, then
# drop the _GoBack bookmark right after it (between "E" and "NTS...").
$splitPoint = $d.Range($start + 34, $start + 34)
$splitPoint.Bookmarks.Add("ZZZ_TEMP_SPLIT")

$goBackPoint = $d.Range($start + 35, $start + 35)
$goBackPoint.Bookmarks.Add("_GoBack")

$tempBm = $d.Bookmarks("ZZZ_TEMP_SPLIT")
$tempBm.Delete()
